$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 230136.83
$ws.Range("J17").Value = 230136.83
$ws.Range("L17").Value = 690410.49
$ws.Range("N17").Value = -690746.49
$ws.Range("H32").Value = 2383.1667
$ws.Range("J32").Value = 2699.8
$ws.Range("L32").Value = 2699.8
$ws.Range("N32").Value = -3351.8
$ws.Range("H53").Value = 247.31818
$ws.Range("I53").Value = 82.84614999999999
$ws.Range("J53").Value = 484.8889
$ws.Range("K53").Value = 82.84614999999999
$ws.Range("L53").Value = 484.8889
$ws.Range("M53").Value = 554.15385
$ws.Range("N53").Value = -1758.8889
$ws.Range("H106").Value = 2533.25
$ws.Range("I106").Value = 2624.875
$ws.Range("K106").Value = 2624.875
$ws.Range("M106").Value = -1993.875
$ws.Range("H129").Value = 942.6053000000001
$ws.Range("I129").Value = 379.9091
$ws.Range("J129").Value = 1171.8518
$ws.Range("K129").Value = 1139.7273
$ws.Range("L129").Value = 3515.5554
$ws.Range("M129").Value = 3860.2727
$ws.Range("N129").Value = -13515.5554
$ws.Range("H137").Value = 1050.0714
$ws.Range("I137").Value = 890.2
$ws.Range("J137").Value = 1449.75
$ws.Range("K137").Value = 2670.6
$ws.Range("L137").Value = 4349.25
$ws.Range("M137").Value = -120.6000000000004
$ws.Range("N137").Value = -9449.25
$ws.Range("H138").Value = 2878.0303
$ws.Range("I138").Value = 615.4054
$ws.Range("J138").Value = 4228.3066
$ws.Range("K138").Value = 1846.2162
$ws.Range("L138").Value = 12684.9198
$ws.Range("M138").Value = 3293.7838
$ws.Range("N138").Value = -22964.9198

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9738.473
$ws.Range("I32").Value = 4299.755
$ws.Range("J32").Value = 54154.668
$ws.Range("K32").Value = 4299.755
$ws.Range("L32").Value = 54154.668
$ws.Range("M32").Value = -4012.755
$ws.Range("N32").Value = -54728.668
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H102").Value = 3581.52
$ws.Range("I102").Value = 2635.9333
$ws.Range("J102").Value = 4999.9
$ws.Range("K102").Value = 2635.9333
$ws.Range("L102").Value = 4999.9
$ws.Range("M102").Value = -1013.9333
$ws.Range("N102").Value = -8243.9
$ws.Range("H122").Value = 1100
$ws.Range("I122").Value = 1084.3334
$ws.Range("J122").Value = 1131.3334
$ws.Range("K122").Value = 3253.0002
$ws.Range("L122").Value = 3394.0002
$ws.Range("M122").Value = -803.0001999999999
$ws.Range("N122").Value = -8294.0002
$ws.Range("H132").Value = 1693.3478
$ws.Range("I132").Value = 702.5714
$ws.Range("J132").Value = 3234.5557
$ws.Range("K132").Value = 2107.7142
$ws.Range("L132").Value = 9703.667099999999
$ws.Range("M132").Value = 422.2857999999997
$ws.Range("N132").Value = -14763.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 111186310
$ws.Range("I20").Value = 109464.5
$ws.Range("J20").Value = 333340000
$ws.Range("K20").Value = 109464.5
$ws.Range("L20").Value = 333340000
$ws.Range("M20").Value = -109217.5
$ws.Range("N20").Value = -333340494
$ws.Range("H99").Value = 55558084
$ws.Range("I99").Value = 76925560
$ws.Range("J99").Value = 2639.4
$ws.Range("K99").Value = 76925560
$ws.Range("L99").Value = 2639.4
$ws.Range("M99").Value = -76924062
$ws.Range("N99").Value = -5635.4
$ws.Range("H134").Value = 76881.92999999999
$ws.Range("I134").Value = 2572.4211
$ws.Range("J134").Value = 253367
$ws.Range("K134").Value = 7717.263300000001
$ws.Range("L134").Value = 760101
$ws.Range("M134").Value = -5182.263300000001
$ws.Range("N134").Value = -765171

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2425.4092
$ws.Range("I31").Value = 2360.0527
$ws.Range("K31").Value = 2360.0527
$ws.Range("M31").Value = -2065.0527
$ws.Range("H34").Value = 2425.4092
$ws.Range("I34").Value = 2360.0527
$ws.Range("K34").Value = 2360.0527
$ws.Range("M34").Value = -2158.0527
$ws.Range("H132").Value = 1501.7667
$ws.Range("I132").Value = 939.6316
$ws.Range("K132").Value = 2818.8948
$ws.Range("M132").Value = -288.8948
$ws.Range("H134").Value = 1958.9773
$ws.Range("I134").Value = 1466.5278
$ws.Range("J134").Value = 4175
$ws.Range("K134").Value = 4399.5834
$ws.Range("L134").Value = 12525
$ws.Range("M134").Value = -1864.5834
$ws.Range("N134").Value = -17595

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1312.4445
$ws.Range("I5").Value = 982.59375
$ws.Range("J5").Value = 3951.25
$ws.Range("K5").Value = 2947.78125
$ws.Range("L5").Value = 11853.75
$ws.Range("M5").Value = -2835.78125
$ws.Range("N5").Value = -12077.75
$ws.Range("H80").Value = 5708.8887
$ws.Range("J80").Value = 5708.8887
$ws.Range("L80").Value = 17126.6661
$ws.Range("N80").Value = -18998.6661
$ws.Range("H83").Value = 5708.8887
$ws.Range("J83").Value = 5708.8887
$ws.Range("L83").Value = 51379.99830000001
$ws.Range("N83").Value = -60739.99830000001
$ws.Range("H127").Value = 2400
$ws.Range("J127").Value = 2400
$ws.Range("L127").Value = 7200
$ws.Range("N127").Value = -17120
$ws.Range("H131").Value = 21629.74
$ws.Range("J131").Value = 2188.8235
$ws.Range("L131").Value = 6566.470499999999
$ws.Range("N131").Value = -16646.4705
$ws.Range("H135").Value = 1312.4445
$ws.Range("I135").Value = 982.59375
$ws.Range("J135").Value = 3951.25
$ws.Range("K135").Value = 8843.34375
$ws.Range("L135").Value = 35561.25
$ws.Range("M135").Value = -6308.34375
$ws.Range("N135").Value = -40631.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4940.8
$ws.Range("I70").Value = 4489.778
$ws.Range("J70").Value = 9000
$ws.Range("K70").Value = 4489.778
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -4219.778
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 4940.8
$ws.Range("I73").Value = 4489.778
$ws.Range("J73").Value = 9000
$ws.Range("K73").Value = 4489.778
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -3553.778
$ws.Range("N73").Value = -10872
$ws.Range("H132").Value = 3230.1738
$ws.Range("I132").Value = 3053.0667
$ws.Range("J132").Value = 3562.25
$ws.Range("K132").Value = 9159.2001
$ws.Range("L132").Value = 10686.75
$ws.Range("M132").Value = -6629.2001
$ws.Range("N132").Value = -15746.75
$ws.Range("H136").Value = 22680
$ws.Range("J136").Value = 22680
$ws.Range("L136").Value = 68040
$ws.Range("N136").Value = -73140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 15875002
$ws.Range("I100").Value = 18520502
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 18520502
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -18519961
$ws.Range("N100").Value = -3082
$ws.Range("H132").Value = 2215.8235
$ws.Range("I132").Value = 1271.9032
$ws.Range("J132").Value = 3678.9
$ws.Range("K132").Value = 3815.7096
$ws.Range("L132").Value = 11036.7
$ws.Range("M132").Value = -1285.7096
$ws.Range("N132").Value = -16096.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2116.3142
$ws.Range("I136").Value = 2105.1875
$ws.Range("J136").Value = 2235
$ws.Range("K136").Value = 6315.5625
$ws.Range("L136").Value = 6705
$ws.Range("M136").Value = -3765.5625
$ws.Range("N136").Value = -11805

